# Applies the edits described by the commit:
#  - rename "Sheet1" -> "Data"
#  - move the saved selection/active cell from A1 to A2 (row 1 is frozen as
#    a header row, so the cursor rests on the first data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Data"

$ws.Range("A2").Select()
